$d = $word.ActiveDocument

# Locate the paragraph whose text is "Listagens:" followed by two tab
# characters (the heading run that needs to be split so the trailing tabs
# lose the underline formatting that "Listagens:" keeps).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Listagens:`t`t`r") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    Write-Host "Target paragraph not found!"
} else {
    $r = $target.Range
    # The run covers "Listagens:" (10 chars) followed by two tabs (2 chars);
    # replace just that run content (not the paragraph mark) with two runs:
    # the original underlined "Listagens:" run, and a new, non-underlined
    # run holding the two tab characters.
    $wholeRange = $d.Range($r.Start, $r.Start + 12)

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
      '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
      '<w:body><w:p><w:r>' +
      '<w:rPr>' +
      '<w:rFonts w:ascii="Roboto Black" w:hAnsi="Roboto Black" w:cs="Roboto Black" w:eastAsia="Roboto Black"/>' +
      '<w:color w:val="auto"/>' +
      '<w:spacing w:val="0"/>' +
      '<w:position w:val="0"/>' +
      '<w:sz w:val="28"/>' +
      '<w:u w:val="single"/>' +
      '<w:shd w:fill="auto" w:val="clear"/>' +
      '</w:rPr>' +
      '<w:t xml:space="preserve">Listagens:</w:t>' +
      '</w:r><w:r>' +
      '<w:rPr>' +
      '<w:rFonts w:ascii="Roboto Black" w:hAnsi="Roboto Black" w:cs="Roboto Black" w:eastAsia="Roboto Black"/>' +
      '<w:color w:val="auto"/>' +
      '<w:spacing w:val="0"/>' +
      '<w:position w:val="0"/>' +
      '<w:sz w:val="28"/>' +
      '<w:shd w:fill="auto" w:val="clear"/>' +
      '</w:rPr>' +
      '<w:tab/><w:tab/>' +
      '</w:r></w:p></w:body></w:document>' +
      '</pkg:xmlData></pkg:part></pkg:package>'

    $wholeRange.InsertXML($xml)
}
